# ForEachTagTemplate.xlsx - add "groupRows" and "groupCols" example sheets
# demonstrating the new GroupTag / groupDir / collapse forEach attributes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "groupRows" - based on the "VertVert" layout (two stacked blocks)
# ---------------------------------------------------------------------
$vertVert = $wb.Worksheets.Item("VertVert")
$vertVert.Range("A1:E3").Copy()

$limitSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$groupRows = $wb.Worksheets.Add($null, $limitSheet)
$groupRows.Name = "groupRows"

# Merge first, then lay down the cell formatting (fills/borders/fonts)
# from the matching template block so the merged range keeps a uniform
# border instead of picking up "outer edge only" merge borders.
$groupRows.Range("A1:E1").Merge()
$groupRows.Range("A1").PasteSpecial(-4122)
$groupRows.Range("A5:E5").Merge()
$groupRows.Range("A5").PasteSpecial(-4122)

# Block 1 (rows 1-3)
$groupRows.Range("A1").Value = '<jt:forEach items="${divisionsList}" var="division" groupDir="rows">Division: ${division.name}'

$groupRows.Range("A2").Value = "City"
$groupRows.Range("B2").Value = "Name"
$groupRows.Range("C2").Value = "Wins"
$groupRows.Range("D2").Value = "Losses"
$groupRows.Range("E2").Value = "Pct."

$groupRows.Range("A3").Value = '<jt:forEach items="${division.teams}" var="team" groupDir="rows" collapse="${division.name.equals(''Central'')}">${team.city}'
$groupRows.Range("B3").Value = '${team.name}'
$groupRows.Range("C3").Value = '${team.wins}'
$groupRows.Range("D3").Value = '${team.losses}'
$groupRows.Range("E3").Value = '${team.pct}</jt:forEach></jt:forEach>'

# Block 2 (rows 5-7) - identical content, repeated further down the sheet
$groupRows.Range("A5").Value = '<jt:forEach items="${divisionsList}" var="division" groupDir="rows">Division: ${division.name}'

$groupRows.Range("A6").Value = "City"
$groupRows.Range("B6").Value = "Name"
$groupRows.Range("C6").Value = "Wins"
$groupRows.Range("D6").Value = "Losses"
$groupRows.Range("E6").Value = "Pct."

$groupRows.Range("A7").Value = '<jt:forEach items="${division.teams}" var="team" groupDir="rows" collapse="${division.name.equals(''Central'')}">${team.city}'
$groupRows.Range("B7").Value = '${team.name}'
$groupRows.Range("C7").Value = '${team.wins}'
$groupRows.Range("D7").Value = '${team.losses}'
$groupRows.Range("E7").Value = '${team.pct}</jt:forEach></jt:forEach>'

[void]$groupRows.Range("A1:E1").Select()

# ---------------------------------------------------------------------
# Sheet "groupCols" - based on the "HorizHoriz" layout (two side-by-side
# vertical blocks)
# ---------------------------------------------------------------------
$horizHoriz = $wb.Worksheets.Item("HorizHoriz")
$horizHoriz.Range("A1:C5").Copy()

$groupCols = $wb.Worksheets.Add($null, $groupRows)
$groupCols.Name = "groupCols"

$groupCols.Range("A1:A5").Merge()
$groupCols.Range("A1").PasteSpecial(-4122)
$groupCols.Range("E1:E5").Merge()
$groupCols.Range("E1").PasteSpecial(-4122)

# Block 1 (columns A-C)
$groupCols.Range("A1").Value = '<jt:forEach items="${divisionsList}" var="division" copyRight="true" groupDir="cols">Division: ${division.name}'
$groupCols.Range("B1").Value = "City"
$groupCols.Range("C1").Value = '<jt:forEach items="${division.teams}" var="team" copyRight="true" groupDir="cols" collapse="${division.name.equals(''Northwest'')}">${team.city}'

$groupCols.Range("B2").Value = "Name"
$groupCols.Range("C2").Value = '${team.name}'

$groupCols.Range("B3").Value = "Wins"
$groupCols.Range("C3").Value = '${team.wins}'

$groupCols.Range("B4").Value = "Losses"
$groupCols.Range("C4").Value = '${team.losses}'

$groupCols.Range("B5").Value = "Pct."
$groupCols.Range("C5").Value = '${team.pct}</jt:forEach></jt:forEach>'

# Block 2 (columns E-G) - identical content, repeated to the right
$groupCols.Range("E1").Value = '<jt:forEach items="${divisionsList}" var="division" copyRight="true" groupDir="cols">Division: ${division.name}'
$groupCols.Range("F1").Value = "City"
$groupCols.Range("G1").Value = '<jt:forEach items="${division.teams}" var="team" copyRight="true" groupDir="cols" collapse="${division.name.equals(''Northwest'')}">${team.city}'

$groupCols.Range("F2").Value = "Name"
$groupCols.Range("G2").Value = '${team.name}'

$groupCols.Range("F3").Value = "Wins"
$groupCols.Range("G3").Value = '${team.wins}'

$groupCols.Range("F4").Value = "Losses"
$groupCols.Range("G4").Value = '${team.losses}'

$groupCols.Range("F5").Value = "Pct."
$groupCols.Range("G5").Value = '${team.pct}</jt:forEach></jt:forEach>'

[void]$groupCols.Range("A1:A5").Select()

# Restore the original active sheet/tab selection.
$vertVert.Activate()

Write-Output "groupRows and groupCols sheets added"
